# Fix the "Number of Lines" column on the methodNumberOfLines sheet:
# several methods (constructors, lambdas) were incorrectly reported as
# having 0 lines; they should read 1 line instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methodNumberOfLines")

# Keep column C (the "Number of Lines" data, which is stored as text in
# this sheet, e.g. "22", "1", "6", ...) formatted as text so the corrected
# values are written back as text too, instead of being auto-converted to
# numbers.
$ws.Range("C2:C21").NumberFormat = "@"

# Rows whose method had an (incorrectly) recorded size of 0 lines:
#   row 4  -> OrderController.lambda$all$0(...)
#   row 8  -> OrderApplicationTests()
#   row 10 -> OrderServiceImpl()
#   row 17 -> KafkaConfig.lambda$table$1(...)
#   row 18 -> KafkaConfig.lambda$stream$0(...)
#   row 21 -> OrderApplication()
$ws.Range("C4").Value = "1"
$ws.Range("C8").Value = "1"
$ws.Range("C10").Value = "1"
$ws.Range("C17").Value = "1"
$ws.Range("C18").Value = "1"
$ws.Range("C21").Value = "1"
